## ---------------------------------------------------------------------
## conv-china.xlsx update:
##  - Sheet1: replace the old 3-column summary with a new Date/HB/QG/
##    QG-HB "Accumulation" table (rows 4-16) and restyle the header row.
##  - "new" sheet: append one more day of data (row 14).
##  - add a new "yishi" sheet holding the data that used to live on the
##    "new" sheet before it was refreshed.
##  - Sheet1's line chart becomes a 3-series "Accumulation" chart over a
##    date axis with a legend; the other two charts on "new" just grow
##    their ranges by one row.
##  - move/resize the chart on Sheet1 and tidy up sheet views/selection.
## ---------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$wsNew = $wb.Worksheets.Item("new")

## 1. Add the "yishi" worksheet after "new", carrying the data that used
##    to be on "new" (pre-refresh) in shifted columns B/C/D.
$wsY = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsNew)
$wsY.Name = "yishi"

$wsY.Range("C2").Value = "HB new"
$wsY.Range("D2").Value = "HB YISHI"

$yishiRows = @(
    @(2.1, 1921, 2606),
    @(2.2, 2103, 3260),
    @(2.3, 2345, 3182),
    @(2.4, 3156, 1957),
    @(2.5, 2987, $null),
    @(2.6, 2447, 2622),
    @(2.7, 2841, 2073),
    @(2.8, 2147, 2067),
    @(2.9, 2618, 2272),
    @(2.1, 2097, 1814),
    @(2.11, 1638, 1685)
)
$r = 3
foreach ($row in $yishiRows) {
    $wsY.Cells.Item($r, 2).Value = $row[0]
    $wsY.Cells.Item($r, 3).Value = $row[1]
    if ($null -ne $row[2]) {
        $wsY.Cells.Item($r, 4).Value = $row[2]
    }
    $r++
}
$wsY.Range("I6").Formula = "=SUM(C3:C13)"
$wsY.Range("J6").Formula = "=SUM(D3:D13)"

$null = $wsY.Range("J7").Select()

## 2. "new" sheet gains one more day of data (row 14); everything else on
##    it is untouched.
$wsNew.Range("A14").Value = 2.12
$wsNew.Range("B14").Value = 14840
$wsNew.Range("C14").Value = 15152
$wsNew.Range("D14").FormulaR1C1 = "=RC[-1]-RC[-2]"

$null = $wsNew.Range("N21").Select()

## 3. Rebuild Sheet1 as the "Accumulation" table: Date / HB / QG / QG-HB.
$ws1.Range("A3:C11").ClearContents()

$ws1.Range("A3").Value = "Date"
$ws1.Range("B3").Value = "HB Accumulation New"
$ws1.Range("C3").Value = "QG Accumulation New"
$ws1.Range("D3").Value = "QG - HB Accumulation New"
$ws1.Range("B3:D3").WrapText = $true
$ws1.Rows("3").RowHeight = 43.2

$sheet1Rows = @(
    @(43862, 9074, 14380),
    @(43863, 11177, 17205),
    @(43864, 13522, 20438),
    @(43865, 16678, 24324),
    @(43866, 19665, 28018),
    @(43867, 22112, 31161),
    @(43868, 24953, 34546),
    @(43869, 27100, 37198),
    @(43870, 29631, 40171),
    @(43871, 31728, 42638),
    @(43872, 33366, 44653),
    @(43873, 48206, 59804),
    @(43874, $null, $null)
)
$r = 4
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 1).NumberFormat = "d-mmm"
    if ($null -ne $row[1]) {
        $ws1.Cells.Item($r, 2).Value = $row[1]
        $ws1.Cells.Item($r, 2).NumberFormat = "#,##0"
    } else {
        $ws1.Cells.Item($r, 2).NumberFormat = "#,##0"
    }
    if ($null -ne $row[2]) {
        $ws1.Cells.Item($r, 3).Value = $row[2]
        $ws1.Cells.Item($r, 3).NumberFormat = "#,##0"
    } else {
        $ws1.Cells.Item($r, 3).NumberFormat = "#,##0"
    }
    if ($r -le 15) {
        $ws1.Cells.Item($r, 4).FormulaR1C1 = "=RC[-1]-RC[-2]"
        $ws1.Cells.Item($r, 4).NumberFormat = "#,##0"
    }
    $r++
}

$ws1.Columns("B").ColumnWidth = 11.29
$ws1.Columns("C").ColumnWidth = 12.29
$ws1.Columns("D").ColumnWidth = 11.92

$null = $ws1.Range("B19").Select()
